# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting (e.g. trailing zeros,
# grouped-thousands-looking strings) instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "52.298.90"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "2.913.72"
$ws.Range("E3").Value = "  +3.85%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "354.31"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "114.40"
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("D7").Value = "0.557"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "0.624"
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").Value = "39.93"
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("D11").Value = "0.0867"
$ws.Range("E11").Value = "  +3.32%  "
$ws.Range("D12").Value = "0.136"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").Value = "19.85"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").Value = "7.75"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "3.373.92"
$ws.Range("E15").Value = "  +3.98%  "
$ws.Range("D16").Value = "2.906.78"
$ws.Range("E16").Value = "  +3.85%  "
$ws.Range("D17").Value = "0.989"
$ws.Range("E17").Value = "  +4.29%  "
$ws.Range("D18").Value = "52.389.17"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").Value = "3.34"
$ws.Range("E19").Value = "  +2.60%  "
$ws.Range("D20").Value = "7.63"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "14.09"
$ws.Range("E21").Value = "  +4.02%  "
$ws.Range("D22").Value = "0.0₃0979"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").Value = "71.18"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").Value = "269.74"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").Value = "2.81"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D26").Value = "0.181"
$ws.Range("E26").Value = "  +12.24%  "
$ws.Range("D27").Value = "26.85"
$ws.Range("E27").Value = "  +2.60%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "10.69"
$ws.Range("E29").Value = "  +2.80%  "
$ws.Range("D30").Value = "0.104"
$ws.Range("E30").Value = "  +15.55%  "
$ws.Range("D31").Value = "6.81"
$ws.Range("E31").Value = "  +11.21%  "
$ws.Range("D32").Value = "37.63"
$ws.Range("E32").Value = "  -4.17%  "
$ws.Range("D33").Value = "2.28"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("D34").Value = "6.11"
$ws.Range("E34").Value = "  +10.81%  "
$ws.Range("D35").Value = "53.13"
$ws.Range("E35").Value = "  +1.84%  "
$ws.Range("D36").Value = "0.0451"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "3.33"
$ws.Range("E38").Value = "  +4.94%  "
$ws.Range("D39").Value = "18.90"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "2.05"
$ws.Range("E40").Value = "  +1.94%  "
$ws.Range("D41").Value = "2.75"
$ws.Range("E41").Value = "  +9.51%  "
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("D43").Value = "23.18"
$ws.Range("E43").Value = "  +5.78%  "
$ws.Range("D44").Value = "118.62"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").Value = "2.18"
$ws.Range("E45").Value = "  -2.12%  "
$ws.Range("D46").Value = "2.53"
$ws.Range("E46").Value = "  +2.04%  "
$ws.Range("D47").Value = "3.54"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").Value = "2.182.52"
$ws.Range("E48").Value = "  +3.30%  "
$ws.Range("D49").Value = "0.261"
$ws.Range("E49").Value = "  +18.07%  "
$ws.Range("D50").Value = "0.0354"
$ws.Range("E50").Value = "  +12.28%  "
$ws.Range("D51").Value = "0.957"
$ws.Range("E51").Value = "  -3.18%  "

# Restore the default (unstyled) look for the Price column cells, matching the
# original workbook which had no explicit style on these cells.
$ws.Range("D2:D51").Style = "Normal"
